$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O9").Value = 2993.77
$ws.Range("N12").Value = 23991.66
$ws.Range("O12").Value = 22805.46
$ws.Range("M16").Value = 70860.57
$ws.Range("K24").Value = 114739.33
